$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$titles = @(
    "PENERAPAN MODEL PEMBELAJARAN DEMONTRASI UNTUK MENINGKATKAN HASIL BELAJAR TIK SISWA KELAS X SMA NEGERI 1 KAWANGKOAN",
    "PENGARUH PENGGUNAAN PEMBELAJARAN E-LEARNING BERBASIS BROWSER BASED TRAINING TERHADAP HASIL PELAJAR KKPI SISWA KELAS X SMK NEGERI I TABUKAN UTARA",
    "PENGEMBANGAN APLIKASI PEMBELAJARAN JARINGAN KOMPUTER BERBASIS ANDROID DI SMK NEGERI 1 AMURANG.",
    "PENERAPAN METODE PEMBELAJARAN DEMONTRASI UNTUK MENINGKATKAN HASIL BELAJAR MERAKIT KOMPUTER SISWA KELAS X TKJ DI SMK N 1 RATAHAN",
    "PENGARUH PEMANFAATAN MEDIA PEMBELAJARAN BERBASIS MULTIMEDIA TERHADAP HASIL BELAJAR TIK SISWA KELAS VIII SMP ADVENT 2 SARIO MANADO"
)

# New header replaces the old "Unnamed: 0" label
$ws.Range("A1").Value = "Judul Skripsi"

# Write the five titles three times in a row, filling rows 2 through 16
$row = 2
for ($copy = 0; $copy -lt 3; $copy++) {
    foreach ($title in $titles) {
        $ws.Cells.Item($row, 1).Value = $title
        $row++
    }
}

# Update the selected range to match the new layout
$ws.Range("A12:A16").Select()

$wb.Save()
